$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.819857
$ws.Range("H2").Value = 5.459571
$ws.Range("I2").Value = 0.01485317462584607
$ws.Range("J2").Value = 0.01485317462584607
$ws.Range("M2").Value = 1.758668
$ws.Range("N2").Value = 5.276004
$ws.Range("O2").Value = 0.02465283256602696
$ws.Range("P2").Value = 0.02465283256602696
$ws.Range("Q2").Value = 3.200524270476
$ws.Range("R2").Value = 28.80471843428401
$ws.Range("S2").Value = 0.0003661728271249432
$ws.Range("T2").Value = 0.0003661728271249433
$ws.Range("G3").Value = 1.819857
$ws.Range("H3").Value = 5.459571
$ws.Range("I3").Value = 0.01485317462584607
$ws.Range("J3").Value = 0.01485317462584607
$ws.Range("O3").Value = 0.0796780206066965
$ws.Range("P3").Value = 0.0796780206066965
$ws.Range("Q3").Value = 10.344102978521
$ws.Range("R3").Value = 93.09692680668901
$ws.Range("S3").Value = 0.001183471553913024
$ws.Range("T3").Value = 0.001183471553913025
$ws.Range("G4").Value = 1.819857
$ws.Range("H4").Value = 5.459571
$ws.Range("I4").Value = 0.01485317462584607
$ws.Range("J4").Value = 0.01485317462584607
$ws.Range("M4").Value = 3.568404
$ws.Range("N4").Value = 10.705212
$ws.Range("O4").Value = 0.05002153126112539
$ws.Range("P4").Value = 0.05002153126112539
$ws.Range("Q4").Value = 6.493984998228
$ws.Range("R4").Value = 58.445864984052
$ws.Range("S4").Value = 0.0007429785388737134
$ws.Range("T4").Value = 0.0007429785388737136
$ws.Range("G5").Value = 1.819857
$ws.Range("H5").Value = 5.459571
$ws.Range("I5").Value = 0.01485317462584607
$ws.Range("J5").Value = 0.01485317462584607
$ws.Range("M5").Value = 60.32626866666666
$ws.Range("N5").Value = 180.978806
$ws.Range("O5").Value = 0.8456476155661511
$ws.Range("P5").Value = 0.8456476155661511
$ws.Range("Q5").Value = 109.785182316914
$ws.Range("R5").Value = 988.066640852226
$ws.Range("S5").Value = 0.01256055170593438
$ws.Range("T5").Value = 0.01256055170593439
$ws.Range("I6").Value = 0.726618572334523
$ws.Range("J6").Value = 0.7266185723345231
$ws.Range("M6").Value = 1.758668
$ws.Range("N6").Value = 5.276004
$ws.Range("O6").Value = 0.02465283256602696
$ws.Range("P6").Value = 0.02465283256602696
$ws.Range("Q6").Value = 156.56992089008
$ws.Range("R6").Value = 1409.12928801072
$ws.Range("S6").Value = 0.01791320600312854
$ws.Range("T6").Value = 0.01791320600312855
$ws.Range("I7").Value = 0.726618572334523
$ws.Range("J7").Value = 0.7266185723345231
$ws.Range("O7").Value = 0.0796780206066965
$ws.Range("P7").Value = 0.0796780206066965
$ws.Range("S7").Value = 0.05789552957967851
$ws.Range("T7").Value = 0.05789552957967852
$ws.Range("I8").Value = 0.726618572334523
$ws.Range("J8").Value = 0.7266185723345231
$ws.Range("M8").Value = 3.568404
$ws.Range("N8").Value = 10.705212
$ws.Range("O8").Value = 0.05002153126112539
$ws.Range("P8").Value = 0.05002153126112539
$ws.Range("Q8").Value = 317.68630121424
$ws.Range("R8").Value = 2859.17671092816
$ws.Range("S8").Value = 0.03634657363094564
$ws.Range("T8").Value = 0.03634657363094565
$ws.Range("I9").Value = 0.726618572334523
$ws.Range("J9").Value = 0.7266185723345231
$ws.Range("M9").Value = 60.32626866666666
$ws.Range("N9").Value = 180.978806
$ws.Range("O9").Value = 0.8456476155661511
$ws.Range("P9").Value = 0.8456476155661511
$ws.Range("Q9").Value = 5370.700503297786
$ws.Range("R9").Value = 48336.30452968008
$ws.Range("S9").Value = 0.6144632631207703
$ws.Range("T9").Value = 0.6144632631207704
$ws.Range("G10").Value = 31.52924033333333
$ws.Range("H10").Value = 94.58772099999999
$ws.Range("I10").Value = 0.257333028084772
$ws.Range("J10").Value = 0.257333028084772
$ws.Range("M10").Value = 1.758668
$ws.Range("N10").Value = 5.276004
$ws.Range("O10").Value = 0.02465283256602696
$ws.Range("P10").Value = 0.02465283256602696
$ws.Range("Q10").Value = 55.44946603854267
$ws.Range("R10").Value = 499.045194346884
$ws.Range("S10").Value = 0.006343988055082599
$ws.Range("T10").Value = 0.0063439880550826
$ws.Range("G11").Value = 31.52924033333333
$ws.Range("H11").Value = 94.58772099999999
$ws.Range("I11").Value = 0.257333028084772
$ws.Range("J11").Value = 0.257333028084772
$ws.Range("O11").Value = 0.0796780206066965
$ws.Range("P11").Value = 0.0796780206066965
$ws.Range("Q11").Value = 179.2128221297265
$ws.Range("R11").Value = 1612.915399167539
$ws.Range("S11").Value = 0.02050378631452207
$ws.Range("T11").Value = 0.02050378631452207
$ws.Range("G12").Value = 31.52924033333333
$ws.Range("H12").Value = 94.58772099999999
$ws.Range("I12").Value = 0.257333028084772
$ws.Range("J12").Value = 0.257333028084772
$ws.Range("M12").Value = 3.568404
$ws.Range("N12").Value = 10.705212
$ws.Range("O12").Value = 0.05002153126112539
$ws.Range("P12").Value = 0.05002153126112539
$ws.Range("Q12").Value = 112.509067322428
$ws.Range("R12").Value = 1012.581605901852
$ws.Range("S12").Value = 0.01287219210886248
$ws.Range("T12").Value = 0.01287219210886248
$ws.Range("G13").Value = 31.52924033333333
$ws.Range("H13").Value = 94.58772099999999
$ws.Range("I13").Value = 0.257333028084772
$ws.Range("J13").Value = 0.257333028084772
$ws.Range("M13").Value = 60.32626866666666
$ws.Range("N13").Value = 180.978806
$ws.Range("O13").Value = 0.8456476155661511
$ws.Range("P13").Value = 0.8456476155661511
$ws.Range("Q13").Value = 1902.041423204569
$ws.Range("R13").Value = 17118.37280884112
$ws.Range("S13").Value = 0.2176130616063049
$ws.Range("T13").Value = 0.2176130616063049
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1464426666666667
$ws.Range("H14").Value = 0.439328
$ws.Range("I14").Value = 0.001195224954858853
$ws.Range("J14").Value = 0.001195224954858853
$ws.Range("M14").Value = 1.758668
$ws.Range("N14").Value = 5.276004
$ws.Range("O14").Value = 0.02465283256602696
$ws.Range("P14").Value = 0.02465283256602696
$ws.Range("Q14").Value = 0.2575440317013333
$ws.Range("R14").Value = 2.317896285312
$ws.Range("S14").Value = 0.00002946568069087243
$ws.Range("T14").Value = 0.00002946568069087243
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1464426666666667
$ws.Range("H15").Value = 0.439328
$ws.Range("I15").Value = 0.001195224954858853
$ws.Range("J15").Value = 0.001195224954858853
$ws.Range("O15").Value = 0.0796780206066965
$ws.Range("P15").Value = 0.0796780206066965
$ws.Range("Q15").Value = 0.8323829973724445
$ws.Range("R15").Value = 7.491446976352
$ws.Range("S15").Value = 0.00009523315858288154
$ws.Range("T15").Value = 0.00009523315858288157
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1464426666666667
$ws.Range("H16").Value = 0.439328
$ws.Range("I16").Value = 0.001195224954858853
$ws.Range("J16").Value = 0.001195224954858853
$ws.Range("M16").Value = 3.568404
$ws.Range("N16").Value = 10.705212
$ws.Range("O16").Value = 0.05002153126112539
$ws.Range("P16").Value = 0.05002153126112539
$ws.Range("Q16").Value = 0.5225665975039999
$ws.Range("R16").Value = 4.703099377536
$ws.Range("S16").Value = 0.00005978698244354927
$ws.Range("T16").Value = 0.00005978698244354929
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.1464426666666667
$ws.Range("H17").Value = 0.439328
$ws.Range("I17").Value = 0.001195224954858853
$ws.Range("J17").Value = 0.001195224954858853
$ws.Range("M17").Value = 60.32626866666666
$ws.Range("N17").Value = 180.978806
$ws.Range("O17").Value = 0.8456476155661511
$ws.Range("P17").Value = 0.8456476155661511
$ws.Range("Q17").Value = 8.834339653596444
$ws.Range("R17").Value = 79.509056882368
$ws.Range("S17").Value = 0.001010739133141549
$ws.Range("T17").Value = 0.001010739133141549
